$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the old row 33 ("OGR001" block), shifting
# everything from the old row 33 onward down by three rows.
$ws.Range("A32:A34").EntireRow.Insert()

# Fill in the new rule rows (REG032..REG034) describing optional
# relationships for the ERD excerpt. Column A is typed out first for the
# two top rows, then column B, matching how the rows were authored.
$ws.Range("A32").Value = "REG032"
$ws.Range("A33").Value = "REG033"
$ws.Range("B32").Value = "Gość może figurować w bazie jednocześnie nie wynajmując żadnego pokoju"
$ws.Range("B33").Value = "Klient może figurować w bazie, nie dokonując jednocześnie żadnej rezerwacji"
$ws.Range("A34").Value = "REG034"
$ws.Range("B34").Value = "Klient może nie wybrać przez pewien czas żadnej metody płatności"

# Restore the active selection to where the author left off editing.
$ws.Activate()
$ws.Range("B34").Select()
